$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells that would otherwise be auto-converted to numbers
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

# Apply the updated values
$ws.Range('D2').Value = '63.360.53'
$ws.Range('E2').Value = '  -1.28%  '
$ws.Range('D3').Value = '3.093.10'
$ws.Range('E3').Value = '  +0.49%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '555.72'
$ws.Range('E5').Value = '  +0.64%  '
$ws.Range('D6').Value = '137.47'
$ws.Range('E6').Value = '  -3.74%  '
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '3.086.06'
$ws.Range('E8').Value = '  +0.52%  '
$ws.Range('D9').Value = '0.496'
$ws.Range('E9').Value = '  +1.41%  '
$ws.Range('D10').Value = '6.64'
$ws.Range('E10').Value = '  +2.39%  '
$ws.Range('E11').Value = '  +5.31%  '
$ws.Range('D12').Value = '0.455'
$ws.Range('E12').Value = '  +1.16%  '
$ws.Range('D13').Value = '35.12'
$ws.Range('E13').Value = '  -1.38%  '
$ws.Range('E14').Value = '  +0.81%  '
$ws.Range('D15').Value = '3.587.45'
$ws.Range('E15').Value = '  +0.71%  '
$ws.Range('D16').Value = '63.322.55'
$ws.Range('E16').Value = '  -1.34%  '
$ws.Range('E17').Value = '  +0.26%  '
$ws.Range('D18').Value = '3.090.74'
$ws.Range('E18').Value = '  +0.48%  '
$ws.Range('D19').Value = '502.23'
$ws.Range('E19').Value = '  +2.82%  '
$ws.Range('E20').Value = '  +1.46%  '
$ws.Range('D21').Value = '13.56'
$ws.Range('E21').Value = '  +0.03%  '
$ws.Range('D22').Value = '0.707'
$ws.Range('E22').Value = '  +3.86%  '
$ws.Range('D23').Value = '7.29'
$ws.Range('E23').Value = '  +1.57%  '
$ws.Range('D24').Value = '78.12'
$ws.Range('E24').Value = '  +0.71%  '
$ws.Range('D25').Value = '12.35'
$ws.Range('E25').Value = '  -0.10%  '
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.26%  '
$ws.Range('D27').Value = '2.76'
$ws.Range('E27').Value = '  +2.00%  '
$ws.Range('D28').Value = '8.18'
$ws.Range('E28').Value = '  -0.64%  '
$ws.Range('E29').Value = '  -1.71%  '
$ws.Range('D30').Value = '0.998'
$ws.Range('E30').Value = '  -0.23%  '
$ws.Range('D31').Value = '26.26'
$ws.Range('E31').Value = '  +2.32%  '
$ws.Range('E32').Value = '  -3.96%  '
$ws.Range('E33').Value = '  -1.30%  '
$ws.Range('D34').Value = '59.59'
$ws.Range('E34').Value = '  +14.42%  '
$ws.Range('D35').Value = '533.89'
$ws.Range('E35').Value = '  -8.24%  '
$ws.Range('E36').Value = '  +0.22%  '
$ws.Range('D37').Value = '5.15'
$ws.Range('E37').Value = '  -3.51%  '
$ws.Range('D38').Value = '0.0413'
$ws.Range('E38').Value = '  +3.26%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').Value = '0.0794'
$ws.Range('E39').Value = '  +1.00%  '
$ws.Range('B40').Value = 'Maker'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D40').Value = '3.065.23'
$ws.Range('E40').Value = '  +2.27%  '
$ws.Range('D41').Value = '0.120'
$ws.Range('E41').Value = '  +1.57%  '
$ws.Range('E42').Value = '  -0.71%  '
$ws.Range('D43').Value = '2.66'
$ws.Range('E43').Value = '  -5.62%  '
$ws.Range('E44').Value = '  +4.41%  '
$ws.Range('D46').Value = '2.07'
$ws.Range('E46').Value = '  -0.67%  '
$ws.Range('D47').Value = '120.45'
$ws.Range('E47').Value = '  +1.69%  '
$ws.Range('D48').Value = '23.94'
$ws.Range('E48').Value = '  -4.36%  '
$ws.Range('D49').Value = '0.107'
$ws.Range('E49').Value = '  -0.21%  '
$ws.Range('B50').Value = 'PEPE'
$ws.Range('C50').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D50').Value = '0.0₃0497'
$ws.Range('E50').Value = '  -4.93%  '
$ws.Range('B51').Value = 'CoreDAO'
$ws.Range('C51').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('D51').Value = '2.33'
$ws.Range('E51').Value = '  +66.82%  '
